# Gain setting sheet: drop the old "Inert"/"Reactant" rows (rows 2-3) and
# append a new "inert" AMU-search row at the bottom, per:
#   "add AMU search of inert in gain setting"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original "Inert" row (row 2) and "Reactant" row (row 3).
# Deleting row 2 twice shifts everything else up by two rows.
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# Append the new inert AMU-search entry as the new last row (row 19).
$ws.Range("A19").Value = "inert"
$ws.Range("B19").Value = "AMU2"

# Match the saved selection in the target workbook.
$ws.Range("B20").Select() | Out-Null
